$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells to text format so the numeric-looking strings (with
# leading spaces / fixed decimal places) are preserved verbatim instead
# of being auto-converted to numbers by Excel.
$cells = @("B4", "D2", "D3", "D4", "F2", "F3", "F4")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 (intrcpt): Chi2 " 1.542" -> " 1.532", p value "0.1230" -> "0.1254"
$ws.Range("D2").Value = " 1.532"
$ws.Range("F2").Value = "0.1254"

# Row 3 (Migratmigrant): Chi2 " 0.105" -> " 0.103", p value "0.7454" -> "0.7483"
$ws.Range("D3").Value = " 0.103"
$ws.Range("F3").Value = "0.7483"

# Row 4 (Pvalue): Estimate "-0.054" -> "-0.052", Chi2 "-1.463" -> "-1.426", p value "0.1433" -> "0.1539"
$ws.Range("B4").Value = "-0.052"
$ws.Range("D4").Value = "-1.426"
$ws.Range("F4").Value = "0.1539"
